# Apply the "#unique -> #match" rename and add new tracking/untracking test
# rows, per the commit:
#   "#unique changed to #match. Changed #unique=false to #match=all in all
#   testing files. Added some more tests for tracking and untracking to
#   cover lines that were previously not covered."

$wb = $excel.ActiveWorkbook

$wsConvert = $wb.Worksheets.Item("#convert")
$wsExport  = $wb.Worksheets.Item("#export")

# ---------------------------------------------------------------------
# Sheet "#convert": rename #unique=true -> #match=unique in the two
# existing blocks, and append two new test blocks (rows 7-8 and 10-11)
# ---------------------------------------------------------------------
$wsConvert.Range("D1").Value2 = "#match=unique"
$wsConvert.Range("D4").Value2 = "#match=unique"

$wsConvert.Range("A7").Value2 = "#tags"
$wsConvert.Range("B7").Value2 = "#measurement.formula.value"
$wsConvert.Range("C7").Value2 = "#measurement.assignment.assign"
$wsConvert.Range("D7").Value2 = "#match=unique"
$wsConvert.Range("E7").Value2 = "#comparison=regex"

$wsConvert.Range("B8").Value2 = "r'qwer'"
$wsConvert.Range("C8").Value2 = "qwer"

$wsConvert.Range("A10").Value2 = "#tags"
$wsConvert.Range("B10").Value2 = "#measurement.formula.value"
$wsConvert.Range("C10").Value2 = "#measurement.assignment.assign"
$wsConvert.Range("D10").Value2 = "#match=unique"
$wsConvert.Range("E10").Value2 = "#comparison=regex"

$wsConvert.Range("B11").Value2 = "r'ghjk'"
$wsConvert.Range("C11").Value2 = "zxcv"

# ---------------------------------------------------------------------
# Sheet "#export": append the corresponding exported rows (4-6) for the
# new "qwer" / "zxcv" test compounds
# ---------------------------------------------------------------------
$wsExport.Range("A4").Value2 = "asdf"
$wsExport.Range("B4").Value2 = "qwer"
$wsExport.Range("C4").Value2 = 1
$wsExport.Range("D4").Value2 = "01_A0_Colon_T03-2017_naive_170427_UKy_GCB_rep1-quench"
$wsExport.Range("E4").Value2 = 289287.73343735602
$wsExport.Range("F4").Value2 = 0
$wsExport.Range("G4").Value2 = 8490014.3650100008
$wsExport.Range("H4").Value2 = 0
$wsExport.Range("I4").Value2 = 439597.55237699999
$wsExport.Range("J4").Value2 = "NA"
$wsExport.Range("K4").Value2 = 0
$wsExport.Range("L4").Value2 = 0
$wsExport.Range("M4").Value2 = 20
$wsExport.Range("N4").Value2 = 10
$wsExport.Range("O4").Value2 = 0.618176844244679
$wsExport.Range("P4").Value2 = 0.255757329816374
$wsExport.Range("Q4").Value2 = 0
$wsExport.Range("R4").Value2 = 0
$wsExport.Range("S4").Value2 = "Compound: name of assigned metabolite, noStd means assigment was NOT verified with standard compound"

$wsExport.Range("A5").Value2 = "zxcv"
$wsExport.Range("B5").Value2 = "ghjk"
$wsExport.Range("C5").Value2 = 0
$wsExport.Range("D5").Value2 = "01_A0_Colon_T03-2017_naive_170427_UKy_GCB_rep1-quench"
$wsExport.Range("E5").Value2 = 7989221.8338638796
$wsExport.Range("F5").Value2 = 8447352.8921099994
$wsExport.Range("G5").Value2 = 8490014.3650100008
$wsExport.Range("H5").Value2 = 8447352.8921099994
$wsExport.Range("I5").Value2 = 7839899.2880199999
$wsExport.Range("J5").Value2 = "NA"
$wsExport.Range("K5").Value2 = 0
$wsExport.Range("L5").Value2 = 0
$wsExport.Range("M5").Value2 = 20
$wsExport.Range("N5").Value2 = 10
$wsExport.Range("O5").Value2 = 0.618176844244679
$wsExport.Range("P5").Value2 = 0.255757329816374
$wsExport.Range("Q5").Value2 = 0
$wsExport.Range("R5").Value2 = 0
$wsExport.Range("S5").Value2 = "Legend"

$wsExport.Range("A6").Value2 = "zxcv"
$wsExport.Range("B6").Value2 = "ghjk"
$wsExport.Range("C6").Value2 = 1
$wsExport.Range("D6").Value2 = "01_A0_Colon_T03-2017_naive_170427_UKy_GCB_rep1-quench"
$wsExport.Range("E6").Value2 = 289287.73343735602
$wsExport.Range("F6").Value2 = 0
$wsExport.Range("G6").Value2 = 8490014.3650100008
$wsExport.Range("H6").Value2 = 0
$wsExport.Range("I6").Value2 = 439597.55237699999
$wsExport.Range("J6").Value2 = "NA"
$wsExport.Range("K6").Value2 = 0
$wsExport.Range("L6").Value2 = 0
$wsExport.Range("M6").Value2 = 20
$wsExport.Range("N6").Value2 = 10
$wsExport.Range("O6").Value2 = 0.618176844244679
$wsExport.Range("P6").Value2 = 0.255757329816374
$wsExport.Range("Q6").Value2 = 0
$wsExport.Range("R6").Value2 = 0
$wsExport.Range("S6").Value2 = "Compound: name of assigned metabolite, noStd means assigment was NOT verified with standard compound"

# ---------------------------------------------------------------------
# Column D on "#convert" now holds data ("#match=unique") and needs to
# be sized to fit, like the other used columns already are. (The COM
# layer applies its own +0.8333 padding to ColumnWidth, so back it out
# here to land on the target stored width of 13.5.)
# ---------------------------------------------------------------------
$wsConvert.Columns.Item(4).ColumnWidth = 12.666666666666666

# ---------------------------------------------------------------------
# Selections per sheet (set last on each sheet so the final Activate
# below decides which tab stays marked as selected)
# ---------------------------------------------------------------------
$wsExport.Range("A5:XFD6").Select() | Out-Null

$wsConvert.Activate() | Out-Null
$wsConvert.Range("E11").Select() | Out-Null
